$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45133
$ws.Range("J2").Value = 80

# Row 3
$ws.Range("D3").Value = 44832

# Row 4
$ws.Range("D4").Value = 44832

# Row 5
$ws.Range("D5").Value = 44846
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 1300
$ws.Range("M5").Value = 1250
$ws.Range("P5").Value = 417

# Row 6
$ws.Range("D6").Value = 44846
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 1000
$ws.Range("P6").Value = 333

# Row 7
$ws.Range("D7").Value = 45134
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = 2500
$ws.Range("P7").Value = 833

# Row 8
$ws.Range("D8").Value = 45135
$ws.Range("J8").Value = 70
$ws.Range("K8").Value = 2500
$ws.Range("L8").Value = 2500
$ws.Range("M8").Value = 2500
$ws.Range("P8").Value = 833

# Row 9
$ws.Range("D9").Value = 44838
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 1200
$ws.Range("L9").Value = 1300
$ws.Range("M9").Value = 1250
$ws.Range("P9").Value = 417

# Row 10
$ws.Range("D10").Value = 44838
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 1000
$ws.Range("P10").Value = 333
